$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the existing 30 data rows (Title in A, Artist in B) before we
#     touch anything. The title strings had a stray leading space baked into
#     the original shared strings -- that gets trimmed. The artist strings
#     (which had a stray trailing space) are carried over unchanged. ---
$data = @(
    @("Smells Like Teen Spirit", "Nirvana "),
    @("...Baby One More Time", "Britney Spears "),
    @("Lose Yourself", "Eminem "),
    @("Wannabe", "Spice Girls "),
    @("I Want It That Way", "Backstreet Boys "),
    @("Beautiful Day", "U2 "),
    @("Waterfalls", "TLC "),
    @("Vogue", "Madonna "),
    @("Creep", "Radiohead "),
    @("Say My Name", "Destiny's Child "),
    @("Losing My Religion", "R.E.M. "),
    @("Hey Ya!", "OutKast "),
    @("I Will Always Love You", "Whitney Houston "),
    @("Wonderwall", "Oasis "),
    @("Vision of Love", "Mariah Carey "),
    @("Basket Case", "Green Day "),
    @("Don't Speak", "No Doubt "),
    @("California Love", "2Pac "),
    @("Bye Bye Bye", "NSYNC "),
    @("Under the Bridge", "Red Hot Chili Peppers "),
    @("Jenny From the Block", "Jennifer Lopez "),
    @("You Oughta Know", "Alanis Morissette "),
    @("Mo Money Mo Problems", "The Notorious B.I.G. "),
    @("Yellow", "Coldplay "),
    @("No Scrubs", "TLC "),
    @("Genie in a Bottle", "Christina Aguilera "),
    @("All the Small Things", "Blink-182"),
    @("Like a Prayer", "Madonna "),
    @("Always Be My Baby", "Mariah Carey "),
    @("Hot in Here", "Nelly ")
)

# --- Insert a fresh row 1 for the header, pushing the 30 song rows down to
#     rows 2-31. ---
$ws.Rows.Item(1).Insert()

# --- Header row: "Titel" / "Artiest", bold. ---
$ws.Range("A1").Value = "Titel"
$ws.Range("B1").Value = "Artiest"
$ws.Range("A1:B1").Font.Bold = $true

# --- Re-write the data rows with the cleaned (trimmed) text. ---
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# --- Column widths. ---
$ws.Range("A1").ColumnWidth = 31.166666666666668
$ws.Range("B1").ColumnWidth = 42.666666666666664

# --- Print setup: portrait orientation. ---
$ws.PageSetup.Orientation = 1

# --- Selection as left by the author. ---
$ws.Range("E12").Select() | Out-Null
